$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay as plain text (they use
# "." as a thousands separator, e.g. "36.666.23", which Excel would
# otherwise reinterpret as a number). Force the Text number format
# before writing the new value, then clear formatting again so the
# cell keeps its original (unstyled) look.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '36.666.23'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '1.964.47'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '244.74'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').Value = '59.27'
$ws.Range('E7').Value = '  +2.09%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '22.35'
$ws.Range('E12').Value = '  +4.15%  '
$ws.Range('D13').Value = '2.253.17'
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '0.830'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '13.72'
$ws.Range('E15').Value = '  +1.71%  '
$ws.Range('D16').Value = '5.26'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '1.967.20'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = '36.585.18'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '70.02'
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '228.77'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '5.06'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +3.32%  '
$ws.Range('D26').Value = '9.23'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  +12.13%  '
$ws.Range('D28').Value = '160.13'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').Value = '19.36'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '4.28'
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = '2.25'
$ws.Range('E36').Value = '  +6.15%  '
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  +13.36%  '
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').Value = '0.0988'
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').Value = '16.13'
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('D45').Value = '1.359.48'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').Value = '87.89'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = '7.14'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('D50').Value = '2.143.75'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').Value = '43.82'
$ws.Range('E51').Value = '  -3.02%  '

$ws.Range('D2').ClearFormats()
$ws.Range('D3').ClearFormats()
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('D12').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D18').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
